# Update the "Förändrad" (Changed) date column (C) for every data row
# from 2026-02-08 (serial 46061) to 2026-02-09 (serial 46062).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 498 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
